$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 235
$ws.Range("B235").Value = 6836277
$ws.Range("C235").Value = "Romania Liga I"
$ws.Range("D235").Value = 45359.625
$ws.Range("E235").Value = "CFR Cluj"
$ws.Range("F235").Value = "AFC Hermannstadt"
$ws.Range("G235").Value = 1
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = "H"
$ws.Range("J235").Value = 1.7
$ws.Range("K235").Value = 3.4
$ws.Range("L235").Value = 5
$ws.Range("M235").Value = 1.65
$ws.Range("N235").Value = 3.5
$ws.Range("O235").Value = 5.25
$ws.Range("P235").Value = -0.75
$ws.Range("Q235").Value = 1.85
$ws.Range("R235").Value = 2
$ws.Range("S235").Value = 2.25
$ws.Range("T235").Value = 1.875
$ws.Range("U235").Value = 1.975
$ws.Range("V235").Value = 0.6499999999999999
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = -1
$ws.Range("Y235").Value = 0.425
$ws.Range("Z235").Value = -0.5
$ws.Range("AA235").Value = -1
$ws.Range("AB235").Value = 0.9750000000000001

# Row 236
$ws.Range("B236").Value = 6852370
$ws.Range("C236").Value = "Romania Liga I"
$ws.Range("D236").Value = 45359.625
$ws.Range("E236").Value = "Dinamo Bucharest"
$ws.Range("F236").Value = "ACS UTA Batrana Doamna"
$ws.Range("G236").Value = 1
$ws.Range("H236").Value = 0
$ws.Range("I236").Value = "H"
$ws.Range("J236").Value = 2.55
$ws.Range("K236").Value = 2.875
$ws.Range("L236").Value = 3
$ws.Range("M236").Value = 2.375
$ws.Range("N236").Value = 3
$ws.Range("O236").Value = 3.1
$ws.Range("P236").Value = -0.25
$ws.Range("Q236").Value = 2
$ws.Range("R236").Value = 1.85
$ws.Range("S236").Value = 2.25
$ws.Range("T236").Value = 1.975
$ws.Range("U236").Value = 1.875
$ws.Range("V236").Value = 1.375
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = 1
$ws.Range("Z236").Value = -1
$ws.Range("AA236").Value = -1
$ws.Range("AB236").Value = 0.875

# Row 237
$ws.Range("B237").Value = 6870268
$ws.Range("C237").Value = "Romania Liga I"
$ws.Range("D237").Value = 45359.625
$ws.Range("E237").Value = "Petrolul Ploiesti"
$ws.Range("F237").Value = "ACS Sepsi"
$ws.Range("G237").Value = 1
$ws.Range("H237").Value = 2
$ws.Range("I237").Value = "A"
$ws.Range("J237").Value = 2.8
$ws.Range("K237").Value = 3
$ws.Range("L237").Value = 2.55
$ws.Range("M237").Value = 3
$ws.Range("N237").Value = 3.2
$ws.Range("O237").Value = 2.3
$ws.Range("P237").Value = 0.25
$ws.Range("Q237").Value = 1.85
$ws.Range("R237").Value = 2
$ws.Range("S237").Value = 2.25
$ws.Range("T237").Value = 1.875
$ws.Range("U237").Value = 1.975
$ws.Range("V237").Value = -1
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = 1.3
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 1
$ws.Range("AA237").Value = 0.875
$ws.Range("AB237").Value = -1

# Row 238
$ws.Range("B238").Value = 6861095
$ws.Range("C238").Value = "Romania Liga I"
$ws.Range("D238").Value = 45359.625
$ws.Range("E238").Value = "FC Botosani"
$ws.Range("F238").Value = "Farul Constanta"
$ws.Range("G238").Value = 0
$ws.Range("H238").Value = 0
$ws.Range("I238").Value = "D"
$ws.Range("J238").Value = 3.75
$ws.Range("K238").Value = 3.4
$ws.Range("L238").Value = 1.909
$ws.Range("M238").Value = 3.1
$ws.Range("N238").Value = 3
$ws.Range("O238").Value = 2.375
$ws.Range("P238").Value = 0.25
$ws.Range("Q238").Value = 1.775
$ws.Range("R238").Value = 2.1
$ws.Range("S238").Value = 2
$ws.Range("T238").Value = 1.8
$ws.Range("U238").Value = 2.05
$ws.Range("V238").Value = -1
$ws.Range("W238").Value = 2
$ws.Range("X238").Value = -1
$ws.Range("Y238").Value = 0.3875
$ws.Range("Z238").Value = -0.5
$ws.Range("AA238").Value = -1
$ws.Range("AB238").Value = 1.05

# Row 239
$ws.Range("B239").Value = 6865915
$ws.Range("C239").Value = "Romania Liga I"
$ws.Range("D239").Value = 45359.625
$ws.Range("E239").Value = "FC Voluntari"
$ws.Range("F239").Value = "Universitatea Cluj"
$ws.Range("G239").Value = 0
$ws.Range("H239").Value = 0
$ws.Range("I239").Value = "D"
$ws.Range("J239").Value = 3.5
$ws.Range("K239").Value = 3.25
$ws.Range("L239").Value = 2.05
$ws.Range("M239").Value = 3.4
$ws.Range("N239").Value = 3.1
$ws.Range("O239").Value = 2.15
$ws.Range("P239").Value = 0.25
$ws.Range("Q239").Value = 1.975
$ws.Range("R239").Value = 1.875
$ws.Range("S239").Value = 2.25
$ws.Range("T239").Value = 2.05
$ws.Range("U239").Value = 1.75
$ws.Range("V239").Value = -1
$ws.Range("W239").Value = 2.1
$ws.Range("X239").Value = -1
$ws.Range("Y239").Value = 0.4875
$ws.Range("Z239").Value = -0.5
$ws.Range("AA239").Value = -1
$ws.Range("AB239").Value = 0.75

# Row 309
$ws.Range("B309").Value = 8191463
$ws.Range("C309").Value = "Romania Liga I"
$ws.Range("D309").Value = 45424.64583333334
$ws.Range("E309").Value = "Dinamo Bucharest"
$ws.Range("F309").Value = "ACS UTA Batrana Doamna"
$ws.Range("G309").Value = 2
$ws.Range("H309").Value = 0
$ws.Range("I309").Value = "H"
$ws.Range("J309").Value = 1.833
$ws.Range("K309").Value = 3.4
$ws.Range("L309").Value = 3.6
$ws.Range("M309").Value = 1.5
$ws.Range("N309").Value = 4.333
$ws.Range("O309").Value = 5
$ws.Range("P309").Value = -1
$ws.Range("Q309").Value = 1.875
$ws.Range("R309").Value = 1.975
$ws.Range("S309").Value = 3
$ws.Range("T309").Value = 2.025
$ws.Range("U309").Value = 1.825
$ws.Range("V309").Value = 0.5
$ws.Range("W309").Value = -1
$ws.Range("X309").Value = -1
$ws.Range("Y309").Value = 0.875
$ws.Range("Z309").Value = -1
$ws.Range("AA309").Value = -1
$ws.Range("AB309").Value = 0.825

# Row 310
$ws.Range("B310").Value = 8191462
$ws.Range("C310").Value = "Romania Liga I"
$ws.Range("D310").Value = 45424.64583333334
$ws.Range("E310").Value = "CSM Politehnica Iasi"
$ws.Range("F310").Value = "Petrolul Ploiesti"
$ws.Range("G310").Value = 2
$ws.Range("H310").Value = 0
$ws.Range("I310").Value = "H"
$ws.Range("J310").Value = 2.1
$ws.Range("K310").Value = 3.3
$ws.Range("L310").Value = 3.1
$ws.Range("M310").Value = 1.8
$ws.Range("N310").Value = 3.2
$ws.Range("O310").Value = 4.2
$ws.Range("P310").Value = -0.5
$ws.Range("Q310").Value = 1.85
$ws.Range("R310").Value = 2
$ws.Range("S310").Value = 2.25
$ws.Range("T310").Value = 2.025
$ws.Range("U310").Value = 1.825
$ws.Range("V310").Value = 0.8
$ws.Range("W310").Value = -1
$ws.Range("X310").Value = -1
$ws.Range("Y310").Value = 0.8500000000000001
$ws.Range("Z310").Value = -1
$ws.Range("AA310").Value = -0.5
$ws.Range("AB310").Value = 0.4125

# Row 311
$ws.Range("B311").Value = 8191475
$ws.Range("C311").Value = "Romania Liga I"
$ws.Range("D311").Value = 45424.64583333334
$ws.Range("E311").Value = "FC U Craiova 1948"
$ws.Range("F311").Value = "AFC Hermannstadt"
$ws.Range("G311").Value = 1
$ws.Range("H311").Value = 3
$ws.Range("I311").Value = "A"
$ws.Range("J311").Value = 2.625
$ws.Range("K311").Value = 3.3
$ws.Range("L311").Value = 2.45
$ws.Range("M311").Value = 2.05
$ws.Range("N311").Value = 3.5
$ws.Range("O311").Value = 3
$ws.Range("P311").Value = -0.25
$ws.Range("Q311").Value = 1.85
$ws.Range("R311").Value = 2
$ws.Range("S311").Value = 2.25
$ws.Range("T311").Value = 1.825
$ws.Range("U311").Value = 2.025
$ws.Range("V311").Value = -1
$ws.Range("W311").Value = -1
$ws.Range("X311").Value = 2
$ws.Range("Y311").Value = -1
$ws.Range("Z311").Value = 1
$ws.Range("AA311").Value = 0.825
$ws.Range("AB311").Value = -1

# Row 312
$ws.Range("B312").Value = 8191523
$ws.Range("C312").Value = "Romania Liga I"
$ws.Range("D312").Value = 45424.64583333334
$ws.Range("E312").Value = "Otelul Galati"
$ws.Range("F312").Value = "FC Botosani"
$ws.Range("G312").Value = 2
$ws.Range("H312").Value = 0
$ws.Range("I312").Value = "H"
$ws.Range("J312").Value = 1.666
$ws.Range("K312").Value = 3.6
$ws.Range("L312").Value = 4.6
$ws.Range("M312").Value = 2.9
$ws.Range("N312").Value = 3.5
$ws.Range("O312").Value = 2.2
$ws.Range("P312").Value = 0.25
$ws.Range("Q312").Value = 1.85
$ws.Range("R312").Value = 2
$ws.Range("S312").Value = 2.25
$ws.Range("T312").Value = 1.875
$ws.Range("U312").Value = 1.975
$ws.Range("V312").Value = 1.9
$ws.Range("W312").Value = -1
$ws.Range("X312").Value = -1
$ws.Range("Y312").Value = 0.8500000000000001
$ws.Range("Z312").Value = -1
$ws.Range("AA312").Value = -0.5
$ws.Range("AB312").Value = 0.4875

# Row 313
$ws.Range("B313").Value = 8191476
$ws.Range("C313").Value = "Romania Liga I"
$ws.Range("D313").Value = 45424.64583333334
$ws.Range("E313").Value = "FC Voluntari"
$ws.Range("F313").Value = "Universitatea Cluj"
$ws.Range("G313").Value = 0
$ws.Range("H313").Value = 1
$ws.Range("I313").Value = "A"
$ws.Range("J313").Value = 3.05
$ws.Range("K313").Value = 3.3
$ws.Range("L313").Value = 2.15
$ws.Range("M313").Value = 2.6
$ws.Range("N313").Value = 3.4
$ws.Range("O313").Value = 2.4
$ws.Range("P313").Value = 0
$ws.Range("Q313").Value = 2
$ws.Range("R313").Value = 1.85
$ws.Range("S313").Value = 2.25
$ws.Range("T313").Value = 2
$ws.Range("U313").Value = 1.85
$ws.Range("V313").Value = -1
$ws.Range("W313").Value = -1
$ws.Range("X313").Value = 1.4
$ws.Range("Y313").Value = -1
$ws.Range("Z313").Value = 0.8500000000000001
$ws.Range("AA313").Value = -1
$ws.Range("AB313").Value = 0.8500000000000001

# Row 319
$ws.Range("Q319").Value = 1.825
$ws.Range("R319").Value = 2.025

Write-Output "done"